# New abbreviations for meter
# Remove three obsolete battery datapoint rows from the "Worksheet" sheet:
#   B_FCE_COUNT (row 24), B_INV_COUNT (row 26), B_LOGIC_BAT_COUNT (row 36)
# Deleting from the bottom up keeps the remaining row numbers stable while
# each delete executes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Rows.Item(36).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(24).Delete()

# Column D ("Description") no longer contains the longest removed string
# (B_LOGIC_BAT_COUNT's description), so its best-fit width shrinks along
# with the row deletions above.
$ws.Columns.Item(4).ColumnWidth = 57.917
